# fix(allowance): issue on import allowance generate payroll
#
# - Fix "Employee id" header casing -> "Employee Id"
# - Apply a Text ("@") number format to the Employee Id column (header cell)
# - Move the active cell / selection to I12
# - Set the sheet's page orientation to portrait

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct the "Employee id" header text to "Employee Id"
$ws.Range("A1").Value = "Employee Id"

# 2) Apply a text number format ("@") to column A's header cell, which is
#    what introduces the extra cellXfs entry (numFmtId 49) used on A1.
$ws.Range("A1").NumberFormat = "@"

# 3) Update the page setup to portrait orientation.
$ps = $ws.PageSetup
$ps.Orientation = 1

# 4) Move the sheet's selection/active cell to I12.
[void]$ws.Range("I12").Select()
